# V 0.55-B55 change: add two new HudBar items to the Tabelle2 "engine merge"
# matrix: AP_VNAV (AP VNAV button) and GPS_TOD (Top-Of-Descend for VNAV
# support). These are inserted as two brand-new columns right before the
# existing "END_OF_COL" marker column, pushing it (and the aircraft-name
# formula column after it) two columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Columns EJ (140) / EK (141) are currently "END_OF_COL" / aircraft-name
# formula. Insert two fresh blank columns there; everything from the old EJ
# onward (including column widths, the END_OF_COL marker and the formula
# column) shifts right to EL/EM automatically, inheriting the existing
# formatting of their left neighbour.
$ws.Columns("EJ:EK").Insert()

# Header row (row 1): new item names.
$ws.Range("EJ1").Value = "AP_VNAV"
$ws.Range("EK1").Value = "GPS_TOD"

# Data rows (2..40): same "|" separator value used by every other item
# column in this matrix.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 140).Value = "|"
    $ws.Cells.Item($r, 141).Value = "|"
}

# Restore the author's view state (selected cell) after the edit.
$ws.Range("EE38").Select() | Out-Null
